# Update the "Förändrad" (Changed) date column (C) for all data rows
# from 2023-09-21 (serial 45190) to 2023-09-23 (serial 45192).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C472").Value = 45192
